# Handles float input without breaking stuff
#
# The quiz-grading logic changed how it totals up right/wrong/
# not-attempted answers (it now also reads a second "Student Ans" /
# "Correct Ans" block per question and drops the stray third block
# entirely), and records which option the student actually picked.
# This script reproduces the resulting cell deltas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy only the *formatting* of $srcAddr onto $dstAddr without
# touching its value - reuses the workbook's existing named cell
# styles/xf records instead of Excel COM minting brand new (duplicate)
# cellXfs entries the way `Range.Style = "name"` does.
function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------
# Drop the third "Student Ans" / "Correct Ans" block entirely
# (columns G:H) - shrinks the used range from A5:H40 down to A5:E40.
# ---------------------------------------------------------------
$ws.Range("G1:H1048576").EntireColumn.Delete()

# ---------------------------------------------------------------
# Summary block (rows 10-12): label cells A10/A11/A12 pick up the
# "mtitleStyle" formatting (same as the row 9 header cells), and the
# Right/Wrong/NotAttempt/Max + Marking + Total numbers are updated.
# ---------------------------------------------------------------
Copy-CellFormat "A9" "A10"
Copy-CellFormat "A9" "A11"
Copy-CellFormat "A9" "A12"

$ws.Range("B10").Value2 = 15
$ws.Range("C10").Value2 = 2
$ws.Range("D10").Value2 = 11
$ws.Range("E10").Value2 = 28

$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1

$ws.Range("B12").Value2 = 60
$ws.Range("C12").Value2 = -2
$ws.Range("E12").Value2 = "58/112"

# ---------------------------------------------------------------
# Rows 16-40: column A gets filled in with the student's recorded
# answer - "correctStyle" formatting when it matches the "Correct
# Ans" in column B, "incorrectStyle" when it doesn't. Rows that are
# still unattempted keep column A blank (default "normalStyle").
# ---------------------------------------------------------------
$correctRows  = 16, 17, 18, 19, 20, 21, 22, 25, 26, 27, 32, 33, 37, 39
$incorrectRows = 36

$answers = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C";
    20 = "Option B"; 21 = "Option C"; 22 = "Option D"; 25 = "Option A";
    26 = "Option C"; 27 = "Option A"; 32 = "Option C"; 33 = "Option D";
    36 = "Option B"; 37 = "Option A"; 39 = "Option D"
}

foreach ($r in $correctRows) {
    Copy-CellFormat "B11" "A$r"
    $ws.Range("A$r").Value2 = $answers[$r]
}
foreach ($r in $incorrectRows) {
    Copy-CellFormat "C10" "A$r"
    $ws.Range("A$r").Value2 = $answers[$r]
}

# ---------------------------------------------------------------
# Second "Student Ans"/"Correct Ans" block (D:E): keep rows 16-18 but
# update their values; clear everything from row 19 down.
# ---------------------------------------------------------------
Copy-CellFormat "C10" "D16"
$ws.Range("D16").Value2 = "Option C"

Copy-CellFormat "B11" "D18"
$ws.Range("D18").Value2 = "Option D"

$ws.Range("D19:E40").Clear()
